# Update EURUSD EDA report with the latest model run values
# (Se actualiza el modelo con la version mas reciente)

$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 43875
$wsSummary.Range("C2").Value = 45982

# --- Sheet: Descriptive Stats ---
$wsDesc = $wb.Worksheets.Item("Descriptive Stats")
$wsDesc.Range("B2").Value = 1.111003366666667
$wsDesc.Range("C2").Value = 0.05783336521899728
$wsDesc.Range("F2").Value = 1.09718
$wsDesc.Range("J2").Value = 0.05205507647786362

# --- Sheet: Returns Stats ---
$wsRet = $wb.Worksheets.Item("Returns Stats")
$wsRet.Range("A2").Value = (5.250214882616561 * [Math]::Pow(10, -5))
$wsRet.Range("B2").Value = 0.004810798658522763
$wsRet.Range("C2").Value = (4.094322411008945 * [Math]::Pow(10, -5))
$wsRet.Range("D2").Value = 0.004808886319010163
$wsRet.Range("E2").Value = 0.1367800467717637
$wsRet.Range("F2").Value = 1.636728195830338
$wsRet.Range("I2").Value = 0.1732447840166188
$wsRet.Range("J2").Value = 170.0556606513871
$wsRet.Range("K2").Value = (1.182721471890143 * [Math]::Pow(10, -37))

# --- Sheet: Stationarity ---
$wsStat = $wb.Worksheets.Item("Stationarity")
$wsStat.Range("B2").Value = -1.723240175213792
$wsStat.Range("C2").Value = 0.4191799631434206
$wsStat.Range("B3").Value = 1.408118139347521
